# Apply the edit described by the commit: "finaler Commit Bericht und Handbuch"
#
# Summary of the underlying data change (Tabelle1 / sheet1):
#   - F22 was empty; a value of 32 is entered.
#   - This cascades through the existing formulas:
#       G22  =F22-E22            -> 14
#       F25  =SUM(F2:F24)        -> 382
#       G25  =SUM(G2:G24)        ->  22
#       F26  =AVERAGE(F2:F24)    -> 18.19047619047619
#   Excel's automatic recalculation takes care of the formula cells once the
#   input cell is written, exactly as it would in the real application.
#
#   - The sheet's on-screen view also changed: the top-left visible cell
#     moved from D1 to B1, and the active/selected cell moved from F26 to
#     F25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Make sure we are working on the right sheet.
$ws.Activate()

# --- Data edit ----------------------------------------------------------
# Enter the new planning value; dependent formulas (G22, F25, G25, F26)
# recalc automatically afterwards.
$ws.Range("F22").Value = 32

# --- View / selection edit ----------------------------------------------
# Scroll the window so column B is the left-most visible column (was D),
# then move the active selection to F25 (was F26).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F25").Select()
